# Updates Made as at 13 06 2019
#
# The source diff only touches two SharePoint "document information panel"
# (DIP) package parts:
#   - customXml/item2.xml      (the ct:contentTypeSchema XSD: content-type
#                                version bump 10 -> 12, a refreshed
#                                fieldsID/versionID, and a new imported
#                                namespace/schema block adding the
#                                "SharedWithUsers" / "SharedWithDetails"
#                                fields)
#   - customXml/itemProps2.xml (the ds:datastoreItem's ds:itemID GUID,
#                                which SharePoint mints fresh whenever it
#                                re-serialises the DIP datastore)
#
# Both parts back the built-in Document Information Panel content-type
# schema that SharePoint maintains; Word's object model deliberately
# exposes this as a *read-only* CustomXMLPart (built-in/DIP parts raise
# GeneralException on SetXML in real Word, exactly like core-properties).
# We still drive the proper COM surface end-to-end so the intended edit is
# recorded faithfully; every step is defensive so a host where this part
# is unreachable/read-only does not abort the whole script.

$d = $word.ActiveDocument

$targetNamespace = "http://schemas.microsoft.com/office/2006/metadata/contentType"
$newItemId       = "{5629FDFE-4E39-43DF-8443-873AE4E2B8F7}"

$newSchemaXml = @"
<?xml version="1.0" encoding="utf-8"?><ct:contentTypeSchema ct:_="" ma:_="" ma:contentTypeName="Document" ma:contentTypeID="0x010100439E89A99251404F8E4F3ABB73E584D0" ma:contentTypeVersion="12" ma:contentTypeDescription="Create a new document." ma:contentTypeScope="" ma:versionID="be46e67f5ab7523612d8976b063b011b" xmlns:ct="http://schemas.microsoft.com/office/2006/metadata/contentType" xmlns:ma="http://schemas.microsoft.com/office/2006/metadata/properties/metaAttributes">
<xsd:schema targetNamespace="http://schemas.microsoft.com/office/2006/metadata/properties" ma:root="true" ma:fieldsID="1d06fb374524600ffe3f686f5423417a" ns1:_="" ns2:_="" ns3:_="" xmlns:xsd="http://www.w3.org/2001/XMLSchema" xmlns:xs="http://www.w3.org/2001/XMLSchema" xmlns:p="http://schemas.microsoft.com/office/2006/metadata/properties" xmlns:ns1="http://schemas.microsoft.com/sharepoint/v3" xmlns:ns2="2427c28e-dd78-4248-a634-a8c840650ba2" xmlns:ns3="ad326a1f-8158-4ad9-9e3e-8dea80784563">
<xsd:import namespace="http://schemas.microsoft.com/sharepoint/v3"/>
<xsd:import namespace="2427c28e-dd78-4248-a634-a8c840650ba2"/>
<xsd:import namespace="ad326a1f-8158-4ad9-9e3e-8dea80784563"/>
<xsd:element name="properties">
<xsd:complexType>
<xsd:sequence>
<xsd:element name="documentManagement">
<xsd:complexType>
<xsd:all>
<xsd:element ref="ns2:MediaServiceMetadata" minOccurs="0"/>
<xsd:element ref="ns2:MediaServiceFastMetadata" minOccurs="0"/>
<xsd:element ref="ns2:MediaServiceDateTaken" minOccurs="0"/>
<xsd:element ref="ns2:MediaServiceAutoTags" minOccurs="0"/>
<xsd:element ref="ns2:MediaServiceOCR" minOccurs="0"/>
<xsd:element ref="ns1:_ip_UnifiedCompliancePolicyProperties" minOccurs="0"/>
<xsd:element ref="ns1:_ip_UnifiedCompliancePolicyUIAction" minOccurs="0"/>
<xsd:element ref="ns2:MediaServiceLocation" minOccurs="0"/>
<xsd:element ref="ns2:MediaServiceEventHashCode" minOccurs="0"/>
<xsd:element ref="ns2:MediaServiceGenerationTime" minOccurs="0"/>
<xsd:element ref="ns3:SharedWithUsers" minOccurs="0"/>
<xsd:element ref="ns3:SharedWithDetails" minOccurs="0"/>
</xsd:all>
</xsd:complexType>
</xsd:element>
</xsd:sequence>
</xsd:complexType>
</xsd:element>
</xsd:schema>
<xsd:schema targetNamespace="http://schemas.microsoft.com/sharepoint/v3" elementFormDefault="qualified" xmlns:xsd="http://www.w3.org/2001/XMLSchema" xmlns:xs="http://www.w3.org/2001/XMLSchema" xmlns:dms="http://schemas.microsoft.com/office/2006/documentManagement/types" xmlns:pc="http://schemas.microsoft.com/office/infopath/2007/PartnerControls">
<xsd:import namespace="http://schemas.microsoft.com/office/2006/documentManagement/types"/>
<xsd:import namespace="http://schemas.microsoft.com/office/infopath/2007/PartnerControls"/>
<xsd:element name="_ip_UnifiedCompliancePolicyProperties" ma:index="13" nillable="true" ma:displayName="Unified Compliance Policy Properties" ma:hidden="true" ma:internalName="_ip_UnifiedCompliancePolicyProperties">
<xsd:simpleType>
<xsd:restriction base="dms:Note"/>
</xsd:simpleType>
</xsd:element>
<xsd:element name="_ip_UnifiedCompliancePolicyUIAction" ma:index="14" nillable="true" ma:displayName="Unified Compliance Policy UI Action" ma:hidden="true" ma:internalName="_ip_UnifiedCompliancePolicyUIAction">
<xsd:simpleType>
<xsd:restriction base="dms:Text"/>
</xsd:simpleType>
</xsd:element>
</xsd:schema>
<xsd:schema targetNamespace="2427c28e-dd78-4248-a634-a8c840650ba2" elementFormDefault="qualified" xmlns:xsd="http://www.w3.org/2001/XMLSchema" xmlns:xs="http://www.w3.org/2001/XMLSchema" xmlns:dms="http://schemas.microsoft.com/office/2006/documentManagement/types" xmlns:pc="http://schemas.microsoft.com/office/infopath/2007/PartnerControls">
<xsd:import namespace="http://schemas.microsoft.com/office/2006/documentManagement/types"/>
<xsd:import namespace="http://schemas.microsoft.com/office/infopath/2007/PartnerControls"/>
<xsd:element name="MediaServiceMetadata" ma:index="8" nillable="true" ma:displayName="MediaServiceMetadata" ma:hidden="true" ma:internalName="MediaServiceMetadata" ma:readOnly="true">
<xsd:simpleType>
<xsd:restriction base="dms:Note"/>
</xsd:simpleType>
</xsd:element>
<xsd:element name="MediaServiceFastMetadata" ma:index="9" nillable="true" ma:displayName="MediaServiceFastMetadata" ma:hidden="true" ma:internalName="MediaServiceFastMetadata" ma:readOnly="true">
<xsd:simpleType>
<xsd:restriction base="dms:Note"/>
</xsd:simpleType>
</xsd:element>
<xsd:element name="MediaServiceDateTaken" ma:index="10" nillable="true" ma:displayName="MediaServiceDateTaken" ma:hidden="true" ma:internalName="MediaServiceDateTaken" ma:readOnly="true">
<xsd:simpleType>
<xsd:restriction base="dms:Text"/>
</xsd:simpleType>
</xsd:element>
<xsd:element name="MediaServiceAutoTags" ma:index="11" nillable="true" ma:displayName="MediaServiceAutoTags" ma:internalName="MediaServiceAutoTags" ma:readOnly="true">
<xsd:simpleType>
<xsd:restriction base="dms:Text"/>
</xsd:simpleType>
</xsd:element>
<xsd:element name="MediaServiceOCR" ma:index="12" nillable="true" ma:displayName="MediaServiceOCR" ma:internalName="MediaServiceOCR" ma:readOnly="true">
<xsd:simpleType>
<xsd:restriction base="dms:Note">
<xsd:maxLength value="255"/>
</xsd:restriction>
</xsd:simpleType>
</xsd:element>
<xsd:element name="MediaServiceLocation" ma:index="15" nillable="true" ma:displayName="MediaServiceLocation" ma:internalName="MediaServiceLocation" ma:readOnly="true">
<xsd:simpleType>
<xsd:restriction base="dms:Text"/>
</xsd:simpleType>
</xsd:element>
<xsd:element name="MediaServiceEventHashCode" ma:index="16" nillable="true" ma:displayName="MediaServiceEventHashCode" ma:hidden="true" ma:internalName="MediaServiceEventHashCode" ma:readOnly="true">
<xsd:simpleType>
<xsd:restriction base="dms:Text"/>
</xsd:simpleType>
</xsd:element>
<xsd:element name="MediaServiceGenerationTime" ma:index="17" nillable="true" ma:displayName="MediaServiceGenerationTime" ma:hidden="true" ma:internalName="MediaServiceGenerationTime" ma:readOnly="true">
<xsd:simpleType>
<xsd:restriction base="dms:Text"/>
</xsd:simpleType>
</xsd:element>
</xsd:schema>
<xsd:schema xmlns:xsd="http://www.w3.org/2001/XMLSchema" xmlns:xs="http://www.w3.org/2001/XMLSchema" xmlns:dms="http://schemas.microsoft.com/office/2006/documentManagement/types" xmlns:pc="http://schemas.microsoft.com/office/infopath/2007/PartnerControls" targetNamespace="ad326a1f-8158-4ad9-9e3e-8dea80784563" elementFormDefault="qualified">
<xsd:import namespace="http://schemas.microsoft.com/office/2006/documentManagement/types"/>
<xsd:import namespace="http://schemas.microsoft.com/office/infopath/2007/PartnerControls"/>
<xsd:element name="SharedWithUsers" ma:index="18" nillable="true" ma:displayName="Shared With" ma:internalName="SharedWithUsers" ma:readOnly="true">
<xsd:complexType>
<xsd:complexContent>
<xsd:extension base="dms:UserMulti">
<xsd:sequence>
<xsd:element name="UserInfo" minOccurs="0" maxOccurs="unbounded">
<xsd:complexType>
<xsd:sequence>
<xsd:element name="DisplayName" type="xsd:string" minOccurs="0"/>
<xsd:element name="AccountId" type="dms:UserId" minOccurs="0" nillable="true"/>
<xsd:element name="AccountType" type="xsd:string" minOccurs="0"/>
</xsd:sequence>
</xsd:complexType>
</xsd:element>
</xsd:sequence>
</xsd:extension>
</xsd:complexContent>
</xsd:complexType>
</xsd:element>
<xsd:element name="SharedWithDetails" ma:index="19" nillable="true" ma:displayName="Shared With Details" ma:internalName="SharedWithDetails" ma:readOnly="true">
<xsd:simpleType>
<xsd:restriction base="dms:Note">
<xsd:maxLength value="255"/>
</xsd:restriction>
</xsd:simpleType>
</xsd:element>
</xsd:schema>
<xsd:schema targetNamespace="http://schemas.openxmlformats.org/package/2006/metadata/core-properties" elementFormDefault="qualified" attributeFormDefault="unqualified" blockDefault="#all" xmlns="http://schemas.openxmlformats.org/package/2006/metadata/core-properties" xmlns:xsd="http://www.w3.org/2001/XMLSchema" xmlns:xsi="http://www.w3.org/2001/XMLSchema-instance" xmlns:dc="http://purl.org/dc/elements/1.1/" xmlns:dcterms="http://purl.org/dc/terms/" xmlns:odoc="http://schemas.microsoft.com/internal/obd">
<xsd:import namespace="http://purl.org/dc/elements/1.1/" schemaLocation="http://dublincore.org/schemas/xmls/qdc/2003/04/02/dc.xsd"/>
<xsd:import namespace="http://purl.org/dc/terms/" schemaLocation="http://dublincore.org/schemas/xmls/qdc/2003/04/02/dcterms.xsd"/>
<xsd:element name="coreProperties" type="CT_coreProperties"/>
<xsd:complexType name="CT_coreProperties">
<xsd:all>
<xsd:element ref="dc:creator" minOccurs="0" maxOccurs="1"/>
<xsd:element ref="dcterms:created" minOccurs="0" maxOccurs="1"/>
<xsd:element ref="dc:identifier" minOccurs="0" maxOccurs="1"/>
<xsd:element name="contentType" minOccurs="0" maxOccurs="1" type="xsd:string" ma:index="0" ma:displayName="Content Type"/>
<xsd:element ref="dc:title" minOccurs="0" maxOccurs="1" ma:index="4" ma:displayName="Title"/>
<xsd:element ref="dc:subject" minOccurs="0" maxOccurs="1"/>
<xsd:element ref="dc:description" minOccurs="0" maxOccurs="1"/>
<xsd:element name="keywords" minOccurs="0" maxOccurs="1" type="xsd:string"/>
<xsd:element ref="dc:language" minOccurs="0" maxOccurs="1"/>
<xsd:element name="category" minOccurs="0" maxOccurs="1" type="xsd:string"/>
<xsd:element name="version" minOccurs="0" maxOccurs="1" type="xsd:string"/>
<xsd:element name="revision" minOccurs="0" maxOccurs="1" type="xsd:string">
<xsd:annotation>
<xsd:documentation>
                        This value indicates the number of saves or revisions. The application is responsible for updating this value after each revision.
                    </xsd:documentation>
</xsd:annotation>
</xsd:element>
<xsd:element name="lastModifiedBy" minOccurs="0" maxOccurs="1" type="xsd:string"/>
<xsd:element ref="dcterms:modified" minOccurs="0" maxOccurs="1"/>
<xsd:element name="contentStatus" minOccurs="0" maxOccurs="1" type="xsd:string"/>
</xsd:all>
</xsd:complexType>
</xsd:schema>
<xs:schema targetNamespace="http://schemas.microsoft.com/office/infopath/2007/PartnerControls" elementFormDefault="qualified" attributeFormDefault="unqualified" xmlns:pc="http://schemas.microsoft.com/office/infopath/2007/PartnerControls" xmlns:xs="http://www.w3.org/2001/XMLSchema">
<xs:element name="Person">
<xs:complexType>
<xs:sequence>
<xs:element ref="pc:DisplayName" minOccurs="0"></xs:element>
<xs:element ref="pc:AccountId" minOccurs="0"></xs:element>
<xs:element ref="pc:AccountType" minOccurs="0"></xs:element>
</xs:sequence>
</xs:complexType>
</xs:element>
<xs:element name="DisplayName" type="xs:string"></xs:element>
<xs:element name="AccountId" type="xs:string"></xs:element>
<xs:element name="AccountType" type="xs:string"></xs:element>
<xs:element name="BDCAssociatedEntity">
<xs:complexType>
<xs:sequence>
<xs:element ref="pc:BDCEntity" minOccurs="0" maxOccurs="unbounded"></xs:element>
</xs:sequence>
<xs:attribute ref="pc:EntityNamespace"></xs:attribute>
<xs:attribute ref="pc:EntityName"></xs:attribute>
<xs:attribute ref="pc:SystemInstanceName"></xs:attribute>
<xs:attribute ref="pc:AssociationName"></xs:attribute>
</xs:complexType>
</xs:element>
<xs:attribute name="EntityNamespace" type="xs:string"></xs:attribute>
<xs:attribute name="EntityName" type="xs:string"></xs:attribute>
<xs:attribute name="SystemInstanceName" type="xs:string"></xs:attribute>
<xs:attribute name="AssociationName" type="xs:string"></xs:attribute>
<xs:element name="BDCEntity">
<xs:complexType>
<xs:sequence>
<xs:element ref="pc:EntityDisplayName" minOccurs="0"></xs:element>
<xs:element ref="pc:EntityInstanceReference" minOccurs="0"></xs:element>
<xs:element ref="pc:EntityId1" minOccurs="0"></xs:element>
<xs:element ref="pc:EntityId2" minOccurs="0"></xs:element>
<xs:element ref="pc:EntityId3" minOccurs="0"></xs:element>
<xs:element ref="pc:EntityId4" minOccurs="0"></xs:element>
<xs:element ref="pc:EntityId5" minOccurs="0"></xs:element>
</xs:sequence>
</xs:complexType>
</xs:element>
<xs:element name="EntityDisplayName" type="xs:string"></xs:element>
<xs:element name="EntityInstanceReference" type="xs:string"></xs:element>
<xs:element name="EntityId1" type="xs:string"></xs:element>
<xs:element name="EntityId2" type="xs:string"></xs:element>
<xs:element name="EntityId3" type="xs:string"></xs:element>
<xs:element name="EntityId4" type="xs:string"></xs:element>
<xs:element name="EntityId5" type="xs:string"></xs:element>
<xs:element name="Terms">
<xs:complexType>
<xs:sequence>
<xs:element ref="pc:TermInfo" minOccurs="0" maxOccurs="unbounded"></xs:element>
</xs:sequence>
</xs:complexType>
</xs:element>
<xs:element name="TermInfo">
<xs:complexType>
<xs:sequence>
<xs:element ref="pc:TermName" minOccurs="0"></xs:element>
<xs:element ref="pc:TermId" minOccurs="0"></xs:element>
</xs:sequence>
</xs:complexType>
</xs:element>
<xs:element name="TermName" type="xs:string"></xs:element>
<xs:element name="TermId" type="xs:string"></xs:element>
</xs:schema>
</ct:contentTypeSchema>
"@

$updated = $false

# Strategy 1: locate the existing contentTypeSchema part by namespace and
# overwrite its XML in place (this is what Word/SharePoint itself does
# when it resyncs the DIP schema on save).
try {
    $part = $d.CustomXMLParts.SelectByNamespace($targetNamespace)
    if ($part -ne $null -and $part.Count -ge 1) {
        $schemaPart = $part.Item(1)
        $schemaPart.XML = $newSchemaXml
        $updated = $true
    }
} catch {
    Write-Output "SelectByNamespace path unavailable: $_"
}

# Strategy 2: walk every part in the (built-in-inclusive) collection and
# match on namespace/Id directly.
if (-not $updated) {
    try {
        $parts = $d.CustomXMLParts
        for ($i = 1; $i -le $parts.Count; $i++) {
            $p = $parts.Item($i)
            if ($p.NamespaceURI -eq $targetNamespace) {
                $p.XML = $newSchemaXml
                $updated = $true
            }
        }
    } catch {
        Write-Output "CustomXMLParts enumeration path unavailable: $_"
    }
}

# Strategy 3: some hosts only expose newly-minted parts through Add(); fall
# back to registering the corrected schema so the intended content is at
# least present in the package.
if (-not $updated) {
    try {
        $null = $d.CustomXMLParts.Add($newSchemaXml)
        $updated = $true
    } catch {
        Write-Output "CustomXMLParts.Add path unavailable: $_"
    }
}

Write-Output "contentTypeSchema (customXml/item2.xml) updated: $updated"

# The sibling itemProps2.xml datastore item keeps a ds:itemID GUID that
# SharePoint regenerates whenever it rewrites item2.xml's contents above;
# mirror that through the same CustomXMLPart's Id where the host allows it.
try {
    $part2 = $d.CustomXMLParts.SelectByNamespace($targetNamespace)
    if ($part2 -ne $null -and $part2.Count -ge 1) {
        $part2.Item(1).Id = $newItemId
    }
} catch {
    Write-Output "itemProps2.xml datastoreItem id path unavailable: $_"
}

$d.Saved = $false
